$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '57.886.03'
$ws.Range('E2').Value = '  -2.62%  '
$ws.Range('D3').Value = '2.282.33'
$ws.Range('E3').Value = '  -2.74%  '
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '529.39'
$ws.Range('E5').Value = '  -5.36%  '
Set-TextValue 'D6' '130.56'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('E7').Value = '  +0.09%  '
Set-TextValue 'D8' '0.582'
$ws.Range('E8').Value = '  +0.72%  '
$ws.Range('D9').Value = '2.282.20'
$ws.Range('E9').Value = '  -2.71%  '
Set-TextValue 'D10' '0.0990'
$ws.Range('E10').Value = '  -4.61%  '
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('E12').Value = '  -0.20%  '
$ws.Range('E13').Value = '  -3.14%  '
Set-TextValue 'D14' '23.34'
$ws.Range('E14').Value = '  -2.80%  '
$ws.Range('D15').Value = '2.690.23'
$ws.Range('E15').Value = '  -2.79%  '
$ws.Range('D16').Value = '57.837.21'
$ws.Range('E16').Value = '  -2.67%  '
$ws.Range('E17').Value = '  -3.40%  '
$ws.Range('D18').Value = '2.295.68'
$ws.Range('E18').Value = '  -2.16%  '
Set-TextValue 'D19' '10.47'
$ws.Range('E19').Value = '  -4.31%  '
Set-TextValue 'D20' '4.15'
$ws.Range('E20').Value = '  -6.18%  '
Set-TextValue 'D21' '310.12'
$ws.Range('E21').Value = '  -2.69%  '
$ws.Range('E22').Value = '  -3.33%  '
$ws.Range('E23').Value = '  -0.19%  '
Set-TextValue 'D24' '62.23'
$ws.Range('E24').Value = '  -2.40%  '
Set-TextValue 'D25' '0.166'
$ws.Range('E25').Value = '  -2.28%  '
Set-TextValue 'D26' '1.00'
$ws.Range('E26').Value = '  +0.05%  '
Set-TextValue 'D27' '7.95'
$ws.Range('E27').Value = '  -4.67%  '
Set-TextValue 'D28' '1.25'
$ws.Range('E28').Value = '  -7.27%  '
Set-TextValue 'D29' '169.53'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  -5.76%  '
$ws.Range('D31').Value = '0.0₃0712'
$ws.Range('E31').Value = '  -4.32%  '
$ws.Range('E32').Value = '  -4.02%  '
$ws.Range('E33').Value = '  -5.65%  '
Set-TextValue 'D34' '0.375'
$ws.Range('E34').Value = '  -5.28%  '
$ws.Range('E35').Value = '  -0.01%  '
Set-TextValue 'D36' '17.70'
$ws.Range('E36').Value = '  -1.59%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  -6.05%  '
$ws.Range('E39').Value = '  -4.85%  '
Set-TextValue 'D40' '38.45'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('E41').Value = '  -5.63%  '
Set-TextValue 'D42' '140.58'
$ws.Range('E42').Value = '  -2.94%  '
Set-TextValue 'D43' '285.42'
$ws.Range('E43').Value = '  -8.38%  '
Set-TextValue 'D44' '3.40'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('E45').Value = '  -0.96%  '
Set-TextValue 'D46' '0.0493'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('E47').Value = '  -2.35%  '
Set-TextValue 'D48' '17.86'
$ws.Range('E48').Value = '  -4.60%  '
$ws.Range('E49').Value = '  -3.31%  '
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('E51').Value = '  -0.78%  '
